# Update "gh-pages output" data refresh for 杭州-漫展信息.xlsx
# Applies updated "想去人数" (want-to-go count) figures and a refreshed
# cover image URL across the 展览 / 演出 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# ---- 展览 (sheet1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 396
$ws1.Range("F5").Value = 34
$ws1.Range("F6").Value = 1255
$ws1.Range("F9").Value = 207
$ws1.Range("F11").Value = 186
$ws1.Range("F12").Value = 1061
$ws1.Range("F15").Value = 201
$ws1.Range("F16").Value = 1534
$ws1.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202405/4W5eYgNl1715222545289.jpeg"
$ws1.Range("F17").Value = 561
$ws1.Range("F19").Value = 360
$ws1.Range("F21").Value = 853
$ws1.Range("F22").Value = 1167
$ws1.Range("F24").Value = 1912
$ws1.Range("F25").Value = 2684
$ws1.Range("F26").Value = 1471
$ws1.Range("F27").Value = 71
$ws1.Range("F28").Value = 51
$ws1.Range("F29").Value = 455
$ws1.Range("F30").Value = 776
$ws1.Range("F31").Value = 1335
$ws1.Range("F32").Value = 836
$ws1.Range("F33").Value = 1420
$ws1.Range("F34").Value = 167
$ws1.Range("F36").Value = 796
$ws1.Range("F37").Value = 651
$ws1.Range("F38").Value = 693
$ws1.Range("F39").Value = 875
$ws1.Range("F41").Value = 262

# ---- 演出 (sheet2) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 655
$ws2.Range("F22").Value = 22

# ---- 全部类型 (sheet4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 396
$ws4.Range("F6").Value = 34
$ws4.Range("F9").Value = 1255
$ws4.Range("F12").Value = 207
$ws4.Range("F14").Value = 186
$ws4.Range("F15").Value = 1061
$ws4.Range("F18").Value = 201
$ws4.Range("F19").Value = 1534
$ws4.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202405/4W5eYgNl1715222545289.jpeg"
$ws4.Range("F20").Value = 561
$ws4.Range("F22").Value = 360
$ws4.Range("F25").Value = 1167
$ws4.Range("F26").Value = 2684
$ws4.Range("F28").Value = 1471
$ws4.Range("F29").Value = 71
$ws4.Range("F31").Value = 51
$ws4.Range("F34").Value = 455
$ws4.Range("F35").Value = 776
$ws4.Range("F36").Value = 1335
$ws4.Range("F39").Value = 836
$ws4.Range("F40").Value = 1420
$ws4.Range("F41").Value = 796
$ws4.Range("F42").Value = 651
$ws4.Range("F43").Value = 693
$ws4.Range("F44").Value = 875
$ws4.Range("F46").Value = 22
$ws4.Range("F48").Value = 262

$wb.Save()
